$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.158.68'
$ws.Range('E2').Value = '  +0.73%  '
$ws.Range('D3').Value = '2.092.65'
$ws.Range('E4').Value = '  +0.04%  '
$style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.17'
$ws.Range('D5').Style = $style
$ws.Range('E5').Value = '  +0.54%  '
$ws.Range('E6').Value = '  +0.55%  '
$style = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.61'
$ws.Range('D7').Style = $style
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +0.11%  '
$style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0846'
$ws.Range('D10').Style = $style
$ws.Range('E10').Value = '  +3.50%  '
$ws.Range('E11').Value = '  +0.43%  '
$ws.Range('D12').Value = '2.403.16'
$ws.Range('E12').Value = '  +2.86%  '
$style = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.61'
$ws.Range('D13').Style = $style
$ws.Range('E13').Value = '  +0.68%  '
$style = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.21'
$ws.Range('D14').Style = $style
$ws.Range('E14').Value = '  +4.15%  '
$ws.Range('E15').Value = '  +6.31%  '
$style = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.774'
$ws.Range('D16').Style = $style
$ws.Range('E16').Value = '  +1.68%  '
$ws.Range('D17').Value = '2.078.67'
$ws.Range('E17').Value = '  +2.43%  '
$ws.Range('D18').Value = '38.105.98'
$ws.Range('E18').Value = '  +0.70%  '
$style = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.01'
$ws.Range('D19').Style = $style
$ws.Range('E19').Value = '  +1.59%  '
$style = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.16'
$ws.Range('D20').Style = $style
$ws.Range('E20').Value = '  +0.63%  '
$ws.Range('D21').Value = '0.0₃0835'
$ws.Range('E21').Value = '  +1.32%  '
$style = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '224.25'
$ws.Range('D22').Style = $style
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('E23').Value = '  +0.00%  '
$style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.40'
$ws.Range('D24').Style = $style
$ws.Range('E24').Value = '  -0.86%  '
$style = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.31'
$ws.Range('D25').Style = $style
$ws.Range('E25').Value = '  +3.11%  '
$style = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '170.08'
$ws.Range('D26').Style = $style
$ws.Range('E26').Value = '  +1.85%  '
$style = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.42'
$ws.Range('D27').Style = $style
$ws.Range('E27').Value = '  +1.21%  '
$style = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.130'
$ws.Range('D28').Style = $style
$ws.Range('E28').Value = '  -0.27%  '
$style = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.97'
$ws.Range('D29').Style = $style
$ws.Range('E29').Value = '  +0.55%  '
$style = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.36'
$ws.Range('D30').Style = $style
$ws.Range('E30').Value = '  +5.90%  '
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('E32').Value = '  +5.74%  '
$ws.Range('E33').Value = '  +4.64%  '
$style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.42'
$ws.Range('D34').Style = $style
$ws.Range('E34').Value = '  +0.34%  '
$style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0604'
$ws.Range('D35').Style = $style
$ws.Range('E35').Value = '  -0.30%  '
$style = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.39'
$ws.Range('D36').Style = $style
$ws.Range('E36').Value = '  +4.71%  '
$ws.Range('E37').Value = '  +1.36%  '
$style = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.54'
$ws.Range('D38').Style = $style
$ws.Range('E38').Value = '  +7.01%  '
$ws.Range('E39').Value = '  -0.08%  '
$style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.99'
$ws.Range('D40').Style = $style
$ws.Range('E40').Value = '  +2.14%  '
$ws.Range('D41').Value = '1.558.24'
$ws.Range('E41').Value = '  +1.85%  '
$style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.95'
$ws.Range('D42').Style = $style
$ws.Range('E42').Value = '  +3.81%  '
$ws.Range('E43').Value = '  +0.54%  '
$style = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.83'
$ws.Range('D44').Style = $style
$ws.Range('E44').Value = '  +1.36%  '
$style = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0908'
$ws.Range('D45').Style = $style
$ws.Range('E45').Value = '  -0.69%  '
$ws.Range('E46').Value = '  +4.92%  '
$style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.11'
$ws.Range('D47').Style = $style
$ws.Range('E47').Value = '  +1.23%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.42'
$ws.Range('D48').Style = $style
$ws.Range('E48').Value = '  +4.91%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$style = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.02'
$ws.Range('D49').Style = $style
$ws.Range('E49').Value = '  +1.64%  '
$ws.Range('D51').Value = '2.290.62'
$ws.Range('E51').Value = '  +2.95%  '
